# Add logic for valid and invalid logins
# Appends four new rows to Sheet1 that exercise combinations of a
# valid/invalid email ("saifzane2@gmail.com" vs "saif") with a
# valid/invalid password ("Boxer@123" vs the numeric 123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: invalid email, "valid" password
$ws.Range("A3").Value = "saif"
$ws.Range("B3").Value = "Boxer@123"

# Row 4: valid email, invalid password
$ws.Range("A4").Value = "saifzane2@gmail.com"
$ws.Range("B4").Value = 123

# Row 5: invalid email, invalid password
$ws.Range("A5").Value = "saif"
$ws.Range("B5").Value = 123

# Row 6: valid email, valid password
$ws.Range("A6").Value = "saifzane2@gmail.com"
$ws.Range("B6").Value = "Boxer@123"
